$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = -13.027
$ws.Range("B3").Value = 5.546
$ws.Range("C3").Value = -12.619
$ws.Range("E3").Value = 16.846
$ws.Range("B4").Value = 6.617
$ws.Range("E5").Value = 16.274
$ws.Range("E6").Value = 16.584
$ws.Range("D8").Value = -8.420000000000002
$ws.Range("C9").Value = -10.988
$ws.Range("E10").Value = 16.802
$ws.Range("A11").Value = -21.645
$ws.Range("D11").Value = -7.604000000000001
$ws.Range("A12").Value = -21.66
$ws.Range("B14").Value = 6.051
$ws.Range("D14").Value = -7.678
$ws.Range("A15").Value = -21.897
$ws.Range("C15").Value = -12.968
$ws.Range("D15").Value = -8.366000000000001
$ws.Range("D17").Value = -8.480999999999998
$ws.Range("C19").Value = -12.869
$ws.Range("C20").Value = -12.173
$ws.Range("E21").Value = 16.771
$ws.Range("C25").Value = -12.549
$ws.Range("B26").Value = 6.325
$ws.Range("D26").Value = -7.800999999999999
$ws.Range("A27").Value = -21.316
$ws.Range("C27").Value = -13.682
$ws.Range("E27").Value = 16.738
$ws.Range("A28").Value = -21.361
$ws.Range("C28").Value = -13.591
$ws.Range("E29").Value = 16.563
$ws.Range("C30").Value = -12.403
$ws.Range("A31").Value = -21.167
$ws.Range("B31").Value = 6.099
$ws.Range("A32").Value = -21.396
$ws.Range("C32").Value = -13.036
$ws.Range("E33").Value = 17.217
$ws.Range("B35").Value = 7.088000000000001
$ws.Range("A36").Value = -20.727
$ws.Range("D36").Value = -8.302
$ws.Range("E36").Value = 16.921
$ws.Range("B37").Value = 6.524000000000001
$ws.Range("A38").Value = -20.369
$ws.Range("B39").Value = 6.603999999999999
$ws.Range("E39").Value = 16.304
$ws.Range("B40").Value = 8.211
$ws.Range("D42").Value = -8.286000000000001
$ws.Range("C44").Value = -13.197
$ws.Range("B45").Value = 5.698
$ws.Range("A46").Value = -21.681
$ws.Range("C47").Value = -12.316
$ws.Range("E47").Value = 16.425
$ws.Range("B52").Value = 4.903
$ws.Range("E53").Value = 17.06
$ws.Range("A54").Value = -21.352
$ws.Range("E54").Value = 16.44
$ws.Range("A55").Value = -21.899
$ws.Range("A56").Value = -21.524
$ws.Range("E56").Value = 16.737
$ws.Range("B57").Value = 5.520999999999999
$ws.Range("C58").Value = -12.75
$ws.Range("E58").Value = 16.635
$ws.Range("E60").Value = 16.405
$ws.Range("C62").Value = -12.975
$ws.Range("D64").Value = -7.784000000000001
$ws.Range("E66").Value = 17.039
$ws.Range("A67").Value = -21.588
$ws.Range("D68").Value = -7.002
$ws.Range("A69").Value = -21.604
$ws.Range("E69").Value = 17.126
$ws.Range("A72").Value = -21.436
$ws.Range("E72").Value = 16.568
$ws.Range("A73").Value = -20.42899999999999
$ws.Range("C77").Value = -13.174
$ws.Range("C78").Value = -12.968
$ws.Range("D79").Value = -7.736
$ws.Range("E80").Value = 16.606
$ws.Range("B81").Value = 6.183999999999999
$ws.Range("E82").Value = 16.65
$ws.Range("A83").Value = -20.506
$ws.Range("B83").Value = 7.114
$ws.Range("E83").Value = 16.618
$ws.Range("C84").Value = -13.292
$ws.Range("C89").Value = -12.188
$ws.Range("D89").Value = -7.424000000000001
$ws.Range("A91").Value = -21.632
$ws.Range("C91").Value = -11.447
$ws.Range("C92").Value = -11.421
$ws.Range("A93").Value = -21.49
$ws.Range("C96").Value = -13.341
$ws.Range("A99").Value = -21.218
$ws.Range("B100").Value = 5.238
$ws.Range("B102").Value = 6.587000000000001
$ws.Range("C102").Value = -12.919
